# CC powerdown implemented, average current ~4mA.
#
# Applies the Tirusse workbook edit:
#  - RX_OFF duration (B25) changed 45 -> 108 (dependent formulas in B27/B28
#    recalc automatically).
#  - "Receive duration" rows (A22:B23) re-styled with the "Calculation"
#    cell style, matching rows 27/28 below them.
#  - Selection moved from D23 to B26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "Calculation" style to A22:B23 (adds s="4" to those cells).
$ws.Range("A22:B23").Style = "Calculation"

# RX_OFF duration: 45 -> 108 (drives B27 and B28 recalculation).
$ws.Range("B25").Value = 108

# Move the active selection to B26.
$ws.Range("B26").Select()
